# Update crypto price/volume data per upstream GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.881.44"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "'2.336.12"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'303.49"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'93.89"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("D10").Value = "'34.00"
$ws.Range("E10").Value = "  -5.25%  "
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").Value = "'18.68"
$ws.Range("E12").Value = "  -4.36%  "
$ws.Range("D13").Value = "'0.121"
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "'6.70"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("D15").Value = "'2.701.79"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "'2.372.90"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").Value = "'0.791"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "'42.834.27"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("E19").Value = "  -5.93%  "
$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "'67.86"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").Value = "'235.52"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("E28").Value = "  -6.75%  "
$ws.Range("D29").Value = "'9.12"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("E30").Value = "  -6.59%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'0.0744"
$ws.Range("E32").Value = "  +5.64%  "
$ws.Range("D33").Value = "'4.97"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").Value = "'17.20"
$ws.Range("E34").Value = "  -4.43%  "
$ws.Range("D35").Value = "'4.36"
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'124.95"
$ws.Range("E38").Value = "  -24.15%  "
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").Value = "'22.35"
$ws.Range("E41").Value = "  +22.80%  "
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").Value = "'1.935.16"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("E45").Value = "  -5.78%  "
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("D47").Value = "'2.70"
$ws.Range("E47").Value = "  -3.84%  "
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").Value = "'2.568.45"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").Value = "'52.76"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").Value = "'71.47"
$ws.Range("E51").Value = "  -1.87%  "
